# Mater update 22 nov 2020
# Refresh branch-wise stock status figures and roll in a new "Branch Total" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ Label = "Branch Total"; Values = @(31, 1, 0, 0, 74) }
    3 = @{ Label = "BOG"; Values = @(35, 14, 5, 14, 38) }
    4 = @{ Label = "BSL"; Values = @(37, 14, 8, 9, 38) }
    5 = @{ Label = "COM"; Values = @(40, 13, 8, 10, 35) }
    6 = @{ Label = "COX"; Values = @(36, 13, 5, 15, 37) }
    7 = @{ Label = "CTG"; Values = @(35, 6, 12, 16, 37) }
    8 = @{ Label = "CTN"; Values = @(32, 14, 4, 10, 46) }
    9 = @{ Label = "DNJ"; Values = @(38, 9, 13, 4, 42) }
    10 = @{ Label = "FEN"; Values = @(33, 6, 5, 12, 50) }
    11 = @{ Label = "FRD"; Values = @(35, 17, 6, 13, 35) }
    12 = @{ Label = "GZP"; Values = @(37, 8, 6, 8, 47) }
    13 = @{ Label = "HZJ"; Values = @(36, 10, 9, 11, 40) }
    14 = @{ Label = "JES"; Values = @(36, 19, 8, 12, 31) }
    15 = @{ Label = "KHL"; Values = @(35, 12, 9, 5, 45) }
    16 = @{ Label = "KRN"; Values = @(30, 11, 5, 7, 53) }
    17 = @{ Label = "KSG"; Values = @(36, 12, 8, 11, 39) }
    18 = @{ Label = "KUS"; Values = @(42, 20, 4, 9, 31) }
    19 = @{ Label = "MHK"; Values = @(32, 6, 3, 8, 57) }
    20 = @{ Label = "MIR"; Values = @(35, 7, 2, 6, 56) }
    21 = @{ Label = "MLV"; Values = @(35, 14, 12, 14, 31) }
    22 = @{ Label = "MOT"; Values = @(34, 11, 6, 3, 52) }
    23 = @{ Label = "MYM"; Values = @(37, 14, 12, 8, 35) }
    24 = @{ Label = "NAJ"; Values = @(44, 15, 8, 12, 27) }
    25 = @{ Label = "NOK"; Values = @(40, 11, 7, 12, 36) }
    26 = @{ Label = "PAT"; Values = @(30, 4, 2, 5, 65) }
    27 = @{ Label = "PBN"; Values = @(32, 15, 7, 9, 43) }
    28 = @{ Label = "RAJ"; Values = @(31, 15, 5, 8, 47) }
    29 = @{ Label = "RNG"; Values = @(36, 12, 5, 8, 45) }
    30 = @{ Label = "SAV"; Values = @(34, 8, 4, 10, 50) }
    31 = @{ Label = "SYL"; Values = @(38, 16, 13, 11, 28) }
    32 = @{ Label = "TGL"; Values = @(38, 15, 7, 10, 36) }
}

foreach ($row in ($data.Keys | Sort-Object)) {
    $info = $data[$row]
    $ws.Cells.Item($row, 1).Value = $info.Label
    $col = 4
    foreach ($val in $info.Values) {
        $ws.Cells.Item($row, $col).Value = $val
        $col = $col + 1
    }
}
